# Adds remaining pages to the backlog
# Fills in the four previously-blank backlog rows (24-27 / sheet rows 28-31)
# on the "Product Backlog" sheet with the new "medication/symptoms" user
# stories, and refreshes the auto-fit row heights across the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Backlog")

# ---------------------------------------------------------------------
# New backlog items (rows 28-31 => IDs 24-27). Values are entered in the
# same order the original author typed them so the shared-strings table
# comes out with the same entries.
# ---------------------------------------------------------------------

# Row 28 / ID 24 - View medication/symptoms page (item + story first)
$ws.Range("C28").Value = "View medication/symptoms page"
$ws.Range("H28").Value = "Given that I am a patient, when I go to view my medication and symptoms, then I have a clear view of the medication that I need to take as well as the ability to edit my symptoms"

# Row 29 / ID 25 - View medication details
$ws.Range("C29").Value = "View medication details"
$ws.Range("H29").Value = "Given that I am a patient, when I click to see details of a medication, then I can see information such as when I last took the medication, when I can next take it, and the script to get it from the pharmacy if needed"

# Row 30 / ID 26 - Edit symptoms
$ws.Range("C30").Value = "Edit symptoms"
$ws.Range("H30").Value = "Given that I am a patient, when I want to edit or add a symptom, then there is a clear way to do so linking from the main view symptoms page to this page"

# Notes column: the common note is entered on row 29 first (then reused via
# copy/paste on rows 30 and 31), and row 28's distinct note is typed after.
$ws.Range("I29").Value = "Accesses and edits user information from the database, will be displayed over the view medication/symptoms page"
$ws.Range("I28").Value = "Accesses and edits user information from the database, there will be some variations between patients and doctors for this page such as doctors linking to adding/editing medication as well as being able to edit the doctors notes on the symptoms"

# Row 31 / ID 27 - Add/edit medication (item + story typed last)
$ws.Range("C31").Value = "Add/edit medication"
$ws.Range("H31").Value = "Given that I am a doctor, when I need to edit or add a users medication, then I can input the details of the medication as well as the script and it will be added to the list of the patients medication"

# Fill in the remaining Owner / Effort / Priority / Status + the two
# Notes cells that reuse row 29's text (no new shared strings introduced).
$ws.Range("D28").Value = "All users"
$ws.Range("E28").Value = 7
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = "To Do "

$ws.Range("D29").Value = "Patient"
$ws.Range("E29").Value = 5
$ws.Range("F29").Value = 4
$ws.Range("G29").Value = "To Do "

$ws.Range("D30").Value = "Patient"
$ws.Range("E30").Value = 7
$ws.Range("F30").Value = 5
$ws.Range("G30").Value = "To Do "
$ws.Range("I30").Value = "Accesses and edits user information from the database, will be displayed over the view medication/symptoms page"

$ws.Range("D31").Value = "Doctor"
$ws.Range("E31").Value = 7
$ws.Range("F31").Value = 6
$ws.Range("G31").Value = "To Do "
$ws.Range("I31").Value = "Accesses and edits user information from the database, will be displayed over the view medication/symptoms page"

# I29 loses its bottom border so it matches the other "Notes" cells in the
# same visual block (copy the no-border format from H30, keeping I29's text).
$ws.Range("H30").Copy() | Out-Null
$ws.Range("I29").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Row heights: the extra content triggers Excel's auto-fit, which (together
# with a slightly different default font-metrics pass) nudges every data
# row's height down a notch and grows a couple of the longer-wrapping rows.
# ---------------------------------------------------------------------

$ws.Range("A5:A7").EntireRow.RowHeight = 51
$ws.Range("A8").EntireRow.RowHeight = 76.5
$ws.Range("A9").EntireRow.RowHeight = 89.25
$ws.Range("A10").EntireRow.RowHeight = 102
$ws.Range("A11").EntireRow.RowHeight = 51
$ws.Range("A12:A13").EntireRow.RowHeight = 63.75
$ws.Range("A14:A17").EntireRow.RowHeight = 51
$ws.Range("A18").EntireRow.RowHeight = 102
$ws.Range("A19:A21").EntireRow.RowHeight = 76.5
$ws.Range("A22:A23").EntireRow.RowHeight = 63.75
$ws.Range("A24").EntireRow.RowHeight = 76.5
$ws.Range("A25").EntireRow.RowHeight = 38.25
$ws.Range("A26:A27").EntireRow.RowHeight = 51
$ws.Range("A28").EntireRow.RowHeight = 63.75
$ws.Range("A29").EntireRow.RowHeight = 76.5
$ws.Range("A30").EntireRow.RowHeight = 51
$ws.Range("A31").EntireRow.RowHeight = 76.5

# ---------------------------------------------------------------------
# View state: scrolled down a bit and the active cell left on the new data
# ---------------------------------------------------------------------

$ws.Activate()
$ws.Range("I32").Select()
